$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H33").Value = 387.3913
$ws.Range("I33").Value = 413.8095
$ws.Range("K33").Value = 413.8095
$ws.Range("M33").Value = -184.8095
$ws.Range("H94").Value = 2658.75
$ws.Range("I94").Value = 2536.818
$ws.Range("K94").Value = 2536.818
$ws.Range("M94").Value = -2085.818
$ws.Range("H125").Value = 1665.4667
$ws.Range("I125").Value = 1360.25
$ws.Range("J125").Value = 2014.2858
$ws.Range("K125").Value = 12242.25
$ws.Range("L125").Value = 18128.5722
$ws.Range("M125").Value = -9782.25
$ws.Range("N125").Value = -23048.5722
$ws.Range("H131").Value = 1311.875
$ws.Range("I131").Value = 865
$ws.Range("K131").Value = 2595
$ws.Range("M131").Value = 2445
$ws.Range("H141").Value = 2268.3076
$ws.Range("I141").Value = 1587.7778
$ws.Range("J141").Value = 3799.5
$ws.Range("K141").Value = 4763.3334
$ws.Range("L141").Value = 11398.5
$ws.Range("M141").Value = 416.6665999999996
$ws.Range("N141").Value = -21758.5

$ws = $wb.Worksheets.Item(2)
$ws.Range("H76").Value = 22000
$ws.Range("J76").Value = 22000
$ws.Range("L76").Value = 22000
$ws.Range("N76").Value = -22676
$ws.Range("H79").Value = 22000
$ws.Range("J79").Value = 22000
$ws.Range("L79").Value = 22000
$ws.Range("N79").Value = -24340
$ws.Range("H122").Value = 2293.9
$ws.Range("I122").Value = 2365.5
$ws.Range("J122").Value = 2007.5
$ws.Range("K122").Value = 7096.5
$ws.Range("L122").Value = 6022.5
$ws.Range("M122").Value = -4646.5
$ws.Range("N122").Value = -10922.5
$ws.Range("H132").Value = 4016.3242
$ws.Range("I132").Value = 3826
$ws.Range("K132").Value = 11478
$ws.Range("M132").Value = -8948

$ws = $wb.Worksheets.Item(3)
$ws.Range("H20").Value = 1567.3214
$ws.Range("I20").Value = 2131.5715
$ws.Range("J20").Value = 1003.0714
$ws.Range("K20").Value = 2131.5715
$ws.Range("L20").Value = 1003.0714
$ws.Range("M20").Value = -1884.5715
$ws.Range("N20").Value = -1497.0714
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H76").Value = 26438
$ws.Range("J76").Value = 26438
$ws.Range("L76").Value = 26438
$ws.Range("N76").Value = -27068
$ws.Range("H79").Value = 26438
$ws.Range("J79").Value = 26438
$ws.Range("L79").Value = 26438
$ws.Range("N79").Value = -28622
$ws.Range("H107").Value = 1022.0769
$ws.Range("I107").Value = 729
$ws.Range("J107").Value = 1999
$ws.Range("K107").Value = 729
$ws.Range("L107").Value = 1999
$ws.Range("M107").Value = 1191
$ws.Range("N107").Value = -5839

$ws = $wb.Worksheets.Item(4)
$ws.Range("H99").Value = 1577.7812
$ws.Range("I99").Value = 1128.0476
$ws.Range("J99").Value = 2436.3635
$ws.Range("K99").Value = 1128.0476
$ws.Range("L99").Value = 2436.3635
$ws.Range("M99").Value = 369.9523999999999
$ws.Range("N99").Value = -5432.363499999999
$ws.Range("H126").Value = 1577.7812
$ws.Range("I126").Value = 1128.0476
$ws.Range("J126").Value = 2436.3635
$ws.Range("K126").Value = 3384.142800000001
$ws.Range("L126").Value = 7309.0905
$ws.Range("M126").Value = -914.1428000000005
$ws.Range("N126").Value = -12249.0905
$ws.Range("H130").Value = 76494.164
$ws.Range("J130").Value = 76494.164
$ws.Range("L130").Value = 76494.164
$ws.Range("N130").Value = -86534.164
$ws.Range("H132").Value = 2241.84
$ws.Range("I132").Value = 1552.4
$ws.Range("J132").Value = 4999.6
$ws.Range("K132").Value = 4657.200000000001
$ws.Range("L132").Value = 14998.8
$ws.Range("M132").Value = -2127.200000000001
$ws.Range("N132").Value = -20058.8

$ws = $wb.Worksheets.Item(5)
$ws.Range("H122").Value = 1650.2727
$ws.Range("I122").Value = 2050.375
$ws.Range("J122").Value = 1421.6428
$ws.Range("K122").Value = 18453.375
$ws.Range("L122").Value = 12794.7852
$ws.Range("M122").Value = -16003.375
$ws.Range("N122").Value = -17694.7852

$ws = $wb.Worksheets.Item(6)
$ws.Range("H102").Value = 1666.05
$ws.Range("I102").Value = 1569.6923
$ws.Range("J102").Value = 1845
$ws.Range("K102").Value = 1569.6923
$ws.Range("L102").Value = 1845
$ws.Range("M102").Value = 52.30770000000007
$ws.Range("N102").Value = -5089
$ws.Range("H107").Value = 588.0741
$ws.Range("I107").Value = 498.25
$ws.Range("J107").Value = 844.7143
$ws.Range("K107").Value = 498.25
$ws.Range("L107").Value = 844.7143
$ws.Range("M107").Value = 1421.75
$ws.Range("N107").Value = -4684.7143
$ws.Range("H113").Value = 20834356
$ws.Range("I113").Value = 62500812
$ws.Range("J113").Value = 1129.25
$ws.Range("K113").Value = 62500812
$ws.Range("L113").Value = 1129.25
$ws.Range("M113").Value = -62498642
$ws.Range("N113").Value = -5469.25
$ws.Range("H122").Value = 1684.25
$ws.Range("I122").Value = 575.2222
$ws.Range("J122").Value = 2591.6365
$ws.Range("K122").Value = 1725.6666
$ws.Range("L122").Value = 7774.9095
$ws.Range("M122").Value = 724.3334
$ws.Range("N122").Value = -12674.9095
$ws.Range("H132").Value = 59672.887
$ws.Range("I132").Value = 101986.95
$ws.Range("J132").Value = 3254.1333
$ws.Range("K132").Value = 305960.85
$ws.Range("L132").Value = 9762.3999
$ws.Range("M132").Value = -303430.85
$ws.Range("N132").Value = -14822.3999

$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 583.1875
$ws.Range("I16").Value = 384.72726
$ws.Range("J16").Value = 1019.8
$ws.Range("K16").Value = 384.72726
$ws.Range("L16").Value = 1019.8
$ws.Range("M16").Value = -214.72726
$ws.Range("N16").Value = -1359.8
$ws.Range("H46").Value = 1639.7307
$ws.Range("I46").Value = 1012.75
$ws.Range("J46").Value = 1918.3889
$ws.Range("K46").Value = 1012.75
$ws.Range("L46").Value = 1918.3889
$ws.Range("M46").Value = -824.75
$ws.Range("N46").Value = -2294.3889
$ws.Range("H61").Value = 2099.5
$ws.Range("I61").Value = 1533.3334
$ws.Range("J61").Value = 2948.75
$ws.Range("K61").Value = 1533.3334
$ws.Range("L61").Value = 2948.75
$ws.Range("M61").Value = -1331.3334
$ws.Range("N61").Value = -3352.75
$ws.Range("H113").Value = 2099.5
$ws.Range("I113").Value = 1533.3334
$ws.Range("J113").Value = 2948.75
$ws.Range("K113").Value = 1533.3334
$ws.Range("L113").Value = 2948.75
$ws.Range("M113").Value = 636.6666
$ws.Range("N113").Value = -7288.75
$ws.Range("H132").Value = 25364
$ws.Range("I132").Value = 42850.418
$ws.Range("J132").Value = 4380.3
$ws.Range("K132").Value = 128551.254
$ws.Range("L132").Value = 13140.9
$ws.Range("M132").Value = -126021.254
$ws.Range("N132").Value = -18200.9

$ws = $wb.Worksheets.Item(8)
$ws.Range("H107").Value = 650
$ws.Range("I107").Value = 650
$ws.Range("K107").Value = 1950
$ws.Range("M107").Value = -30
$ws.Range("H132").Value = 2281.4
$ws.Range("I132").Value = 1724.3
$ws.Range("J132").Value = 3395.6
$ws.Range("K132").Value = 5172.9
$ws.Range("L132").Value = 10186.8
$ws.Range("M132").Value = -2642.9
$ws.Range("N132").Value = -15246.8
